$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.414.68"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "1.592.31"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.69"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.506"
$ws.Range("E6").Value = "  -0.93%  "
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0612"
$ws.Range("E8").Value = "  -0.95%  "
$ws.Range("E9").Value = "  -1.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.46"
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0845"
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("D12").Value = "1.816.77"
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("D13").Value = "1.590.14"
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.03"
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.518"
$ws.Range("E15").Value = "  -1.01%  "
$ws.Range("E16").Value = "  -1.46%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.0₃0734"
$ws.Range("E17").Value = "  -2.67%  "
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.01"
$ws.Range("E18").Value = "  +0.22%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "206.93"
$ws.Range("E19").Value = "  -1.10%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.03"
$ws.Range("E20").Value = "  -2.51%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.26"
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("B22").Value = "Toncoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.21"
$ws.Range("E22").Value = "  -4.45%  "
$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.88"
$ws.Range("E23").Value = "  -0.55%  "
$ws.Range("B24").Value = "Monero"
$ws.Range("C24").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "144.66"
$ws.Range("E24").Value = "  +1.67%  "
$ws.Range("B25").Value = "BinanceUSD"
$ws.Range("C25").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.01"
$ws.Range("E25").Value = "  +0.25%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.06"
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.113"
$ws.Range("E27").Value = "  -1.21%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.26"
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0503"
$ws.Range("E29").Value = "  -3.75%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.14"
$ws.Range("E30").Value = "  -0.85%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.23"
$ws.Range("E31").Value = "  -0.88%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.94"
$ws.Range("E32").Value = "  -1.35%  "
$ws.Range("B33").Value = "Maker"
$ws.Range("C33").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D33").Value = "1.279.23"
$ws.Range("E33").Value = "  -1.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.25"
$ws.Range("E34").Value = "  +13.18%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.49"
$ws.Range("E35").Value = "  +1.03%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.597"
$ws.Range("E36").Value = "  -4.30%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.47"
$ws.Range("E37").Value = "  -2.16%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0167"
$ws.Range("E38").Value = "  -2.30%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.814"
$ws.Range("E39").Value = "  -1.28%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.40"
$ws.Range("E40").Value = "  -0.40%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.21"
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.767"
$ws.Range("E42").Value = "  -1.91%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "62.16"
$ws.Range("E43").Value = "  -1.46%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "1.728.01"
$ws.Range("E44").Value = "  -0.49%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "88.99"
$ws.Range("E45").Value = "  -2.56%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.56"
$ws.Range("E46").Value = "  -0.93%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.102"
$ws.Range("E47").Value = "  +1.10%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0100"
$ws.Range("E48").Value = "  -3.96%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0511"
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.01"
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.43"
$ws.Range("E51").Value = "  +0.62%  "
